$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) "67.655.25"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.19%  "

Set-TextValue $ws.Cells.Item(3, 4) "3.316.16"
Set-TextValue $ws.Cells.Item(3, 5) "  +0.19%  "

Set-TextValue $ws.Cells.Item(4, 4) "1.00"
Set-TextValue $ws.Cells.Item(4, 5) "  +0.17%  "

Set-TextValue $ws.Cells.Item(5, 4) "580.77"
Set-TextValue $ws.Cells.Item(5, 5) "  -0.02%  "

Set-TextValue $ws.Cells.Item(6, 4) "174.69"
Set-TextValue $ws.Cells.Item(6, 5) "  -4.11%  "

Set-TextValue $ws.Cells.Item(7, 5) "  +0.08%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.587"
Set-TextValue $ws.Cells.Item(8, 5) "  -0.33%  "

Set-TextValue $ws.Cells.Item(9, 4) "3.314.62"
Set-TextValue $ws.Cells.Item(9, 5) "  +0.30%  "

Set-TextValue $ws.Cells.Item(10, 4) "0.177"
Set-TextValue $ws.Cells.Item(10, 5) "  -0.15%  "

Set-TextValue $ws.Cells.Item(11, 4) "0.575"
Set-TextValue $ws.Cells.Item(11, 5) "  -0.50%  "

Set-TextValue $ws.Cells.Item(12, 4) "45.30"
Set-TextValue $ws.Cells.Item(12, 5) "  -2.29%  "

Set-TextValue $ws.Cells.Item(13, 4) "0.0000268"
Set-TextValue $ws.Cells.Item(13, 5) "  -1.57%  "

Set-TextValue $ws.Cells.Item(14, 4) "657.31"
Set-TextValue $ws.Cells.Item(14, 5) "  +4.60%  "

Set-TextValue $ws.Cells.Item(15, 4) "3.864.73"
Set-TextValue $ws.Cells.Item(15, 5) "  +0.56%  "

Set-TextValue $ws.Cells.Item(16, 4) "8.37"
Set-TextValue $ws.Cells.Item(16, 5) "  -0.96%  "

Set-TextValue $ws.Cells.Item(17, 4) "67.630.70"
Set-TextValue $ws.Cells.Item(17, 5) "  -0.38%  "

Set-TextValue $ws.Cells.Item(18, 5) "  -0.81%  "

Set-TextValue $ws.Cells.Item(19, 4) "3.327.28"
Set-TextValue $ws.Cells.Item(19, 5) "  +0.29%  "

Set-TextValue $ws.Cells.Item(20, 4) "17.29"
Set-TextValue $ws.Cells.Item(20, 5) "  -2.10%  "

Set-TextValue $ws.Cells.Item(21, 4) "10.94"
Set-TextValue $ws.Cells.Item(21, 5) "  +0.37%  "

Set-TextValue $ws.Cells.Item(22, 4) "0.884"
Set-TextValue $ws.Cells.Item(22, 5) "  -1.82%  "

Set-TextValue $ws.Cells.Item(23, 4) "5.32"
Set-TextValue $ws.Cells.Item(23, 5) "  +4.95%  "

Set-TextValue $ws.Cells.Item(24, 4) "16.91"
Set-TextValue $ws.Cells.Item(24, 5) "  -3.81%  "

Set-TextValue $ws.Cells.Item(25, 4) "98.27"
Set-TextValue $ws.Cells.Item(25, 5) "  +1.65%  "

Set-TextValue $ws.Cells.Item(26, 5) "  -3.36%  "

Set-TextValue $ws.Cells.Item(27, 4) "2.65"
Set-TextValue $ws.Cells.Item(27, 5) "  -4.13%  "

Set-TextValue $ws.Cells.Item(28, 4) "9.22"
Set-TextValue $ws.Cells.Item(28, 5) "  -3.55%  "

Set-TextValue $ws.Cells.Item(29, 4) "33.34"
Set-TextValue $ws.Cells.Item(29, 5) "  +3.14%  "

Set-TextValue $ws.Cells.Item(30, 4) "8.39"
Set-TextValue $ws.Cells.Item(30, 5) "  -2.16%  "

Set-TextValue $ws.Cells.Item(31, 5) "  +6.70%  "

Set-TextValue $ws.Cells.Item(32, 4) "567.58"
Set-TextValue $ws.Cells.Item(32, 5) "  -4.91%  "

Set-TextValue $ws.Cells.Item(33, 4) "10.90"
Set-TextValue $ws.Cells.Item(33, 5) "  -0.37%  "

Set-TextValue $ws.Cells.Item(34, 4) "0.104"
Set-TextValue $ws.Cells.Item(34, 5) "  -0.05%  "

Set-TextValue $ws.Cells.Item(35, 5) "  +0.26%  "

Set-TextValue $ws.Cells.Item(36, 4) "3.661.05"
Set-TextValue $ws.Cells.Item(36, 5) "  -6.92%  "

Set-TextValue $ws.Cells.Item(37, 4) "56.11"
Set-TextValue $ws.Cells.Item(37, 5) "  +0.57%  "

Set-TextValue $ws.Cells.Item(38, 4) "3.26"
Set-TextValue $ws.Cells.Item(38, 5) "  -6.47%  "

Set-TextValue $ws.Cells.Item(39, 4) "34.14"
Set-TextValue $ws.Cells.Item(39, 5) "  +4.67%  "

Set-TextValue $ws.Cells.Item(40, 5) "  +1.38%  "

Set-TextValue $ws.Cells.Item(41, 5) "  -1.99%  "

Set-TextValue $ws.Cells.Item(42, 5) "  -4.18%  "

Set-TextValue $ws.Cells.Item(43, 4) "3.36"
Set-TextValue $ws.Cells.Item(43, 5) "  -1.40%  "

Set-TextValue $ws.Cells.Item(44, 5) "  -1.63%  "

Set-TextValue $ws.Cells.Item(45, 4) "0.0₃0659"
Set-TextValue $ws.Cells.Item(45, 5) "  -3.50%  "

Set-TextValue $ws.Cells.Item(46, 4) "0.0405"
Set-TextValue $ws.Cells.Item(46, 5) "  -1.85%  "

Set-TextValue $ws.Cells.Item(47, 4) "2.59"
Set-TextValue $ws.Cells.Item(47, 5) "  +1.80%  "

Set-TextValue $ws.Cells.Item(48, 4) "0.127"
Set-TextValue $ws.Cells.Item(48, 5) "  -0.87%  "

Set-TextValue $ws.Cells.Item(49, 5) "  -0.26%  "

Set-TextValue $ws.Cells.Item(50, 5) "  -1.79%  "

Set-TextValue $ws.Cells.Item(51, 4) "129.39"
Set-TextValue $ws.Cells.Item(51, 5) "  -0.75%  "
